$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Tabela Livro" box: the "publicacao"/IntegerField field becomes "data_publicacao"/DateField ---
$ws.Range("H7").Value = "data_publicacao"
$ws.Range("I7").Value = "DateField"

# --- "Tabela Categoria" box shrinks: remove the "descricao"/TextField row ---
# B8:C8 ("senha"/CharField) picks up the plain thin-box style already used
# elsewhere in the sheet (e.g. H5:I5)
$ws.Range("H5:I5").Copy() | Out-Null
$ws.Range("B8:C8").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Clear the contents of K5:L5 (the descricao / TextField row)
$ws.Range("K5:L5").ClearContents()
# That row is no longer part of the box, so drop its border entirely
$ws.Range("K5:L5").Borders.LineStyle = -4142

# The header (K3:L3) no longer needs a bottom border since the box now closes
# on the row right below it (K4:L4)
$ws.Range("K3:L3").Borders.Item(9).LineStyle = -4142

# K4:L4 ("nome"/CharField) is now the last row of the box, so it gets a full
# thin border on all sides (matching the style already used elsewhere, e.g. H5:I5)
$ws.Range("H5:I5").Copy() | Out-Null
$ws.Range("K4:L4").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# --- Column H needs to be a bit wider to fit "data_publicacao" ---
$ws.Columns("H").ColumnWidth = 14.6

# --- Restore the selected cell ---
$ws.Range("K8").Select() | Out-Null
